$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1863256666666667
$ws.Range("N2").Value = 0.5589770000000001
$ws.Range("O2").Value = 0.01657678358851065
$ws.Range("P2").Value = 0.01657678358851065
$ws.Range("Q2").Value = 0.03838979505733334
$ws.Range("R2").Value = 0.345508155516
$ws.Range("S2").Value = 0.007683335408215065
$ws.Range("T2").Value = 0.007683335408215065
$ws.Range("O3").Value = 0.5186672939413604
$ws.Range("P3").Value = 0.5186672939413604
$ws.Range("S3").Value = 0.2404021723119321
$ws.Range("T3").Value = 0.2404021723119321
$ws.Range("M4").Value = 5.212463666666667
$ws.Range("N4").Value = 15.637391
$ws.Range("O4").Value = 0.4637358003923671
$ws.Range("P4").Value = 0.4637358003923669
$ws.Range("Q4").Value = 1.073955164025334
$ws.Range("R4").Value = 9.665596476228
$ws.Range("S4").Value = 0.2149414375947554
$ws.Range("T4").Value = 0.2149414375947554
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.01146633333333333
$ws.Range("N5").Value = 0.034399
$ws.Range("O5").Value = 0.001020122077762015
$ws.Range("P5").Value = 0.001020122077762015
$ws.Range("Q5").Value = 0.002362477454666667
$ws.Range("R5").Value = 0.021262297092
$ws.Range("S5").Value = 0.0004728263501131353
$ws.Range("T5").Value = 0.0004728263501131352
$ws.Range("G6").Value = 0.2384863333333333
$ws.Range("H6").Value = 0.715459
$ws.Range("I6").Value = 0.5365002283349842
$ws.Range("J6").Value = 0.5365002283349842
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.1863256666666667
$ws.Range("N6").Value = 0.5589770000000001
$ws.Range("O6").Value = 0.01657678358851065
$ws.Range("P6").Value = 0.01657678358851065
$ws.Range("Q6").Value = 0.04443612504922223
$ws.Range("R6").Value = 0.399925125443
$ws.Range("S6").Value = 0.008893448180295584
$ws.Range("T6").Value = 0.008893448180295582
$ws.Range("G7").Value = 0.2384863333333333
$ws.Range("H7").Value = 0.715459
$ws.Range("I7").Value = 0.5365002283349842
$ws.Range("J7").Value = 0.5365002283349842
$ws.Range("O7").Value = 0.5186672939413604
$ws.Range("P7").Value = 0.5186672939413604
$ws.Range("S7").Value = 0.2782651216294282
$ws.Range("T7").Value = 0.2782651216294282
$ws.Range("G8").Value = 0.2384863333333333
$ws.Range("H8").Value = 0.715459
$ws.Range("I8").Value = 0.5365002283349842
$ws.Range("J8").Value = 0.5365002283349842
$ws.Range("M8").Value = 5.212463666666667
$ws.Range("N8").Value = 15.637391
$ws.Range("O8").Value = 0.4637358003923671
$ws.Range("P8").Value = 0.4637358003923669
$ws.Range("Q8").Value = 1.243101347496556
$ws.Range("R8").Value = 11.187912127469
$ws.Range("S8").Value = 0.2487943627976116
$ws.Range("T8").Value = 0.2487943627976115
$ws.Range("G9").Value = 0.2384863333333333
$ws.Range("H9").Value = 0.715459
$ws.Range("I9").Value = 0.5365002283349842
$ws.Range("J9").Value = 0.5365002283349842
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.01146633333333333
$ws.Range("N9").Value = 0.034399
$ws.Range("O9").Value = 0.001020122077762015
$ws.Range("P9").Value = 0.001020122077762015
$ws.Range("Q9").Value = 0.002734563793444444
$ws.Range("R9").Value = 0.024611074141
$ws.Range("S9").Value = 0.0005472957276488796
$ws.Range("T9").Value = 0.0005472957276488795
